$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell I1 ("url"), styled like the other header cells (copy H1's
# formatting - bold, centered, bordered - onto I1, then set its text).
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("I1").Value = "url"

# URLs for each fonction/categorie row (I2:I11)
$urls = @(
    "https://www.archives13.fr/n/archives-anciennes/n:101",
    "https://www.archives13.fr/n/archives-revolutionnaires/n:102",
    "https://www.archives13.fr/n/archives-modernes-et-contemporaines/n:103",
    "https://www.archives13.fr/n/archives-hospitalieres/n:104",
    "https://www.archives13.fr/n/archives-communales-et-intercommunales-deposees/n:105",
    "https://www.archives13.fr/n/archives-privees/n:106",
    "https://www.archives13.fr/n/fonds-iconographiques-et-audiovisuels/n:107",
    "https://www.archives13.fr/n/bibliotheque/n:108",
    "https://www.archives13.fr/n/etat-civil/n:109",
    "https://www.archives13.fr/n/archives-notariales/n:110"
)

for ($i = 0; $i -lt $urls.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $urls[$i]
}
